# Update the "paises" (countries) COVID-19 dashboard with a newer data
# snapshot: refresh the per-country counters that moved between the
# 15:35 and 16:52 pulls, then re-sort the table (rows 4:219) descending
# by "Casos totales" (column B) so the leaderboard reflects the new
# totals - this is what reshuffles Portugal/Kazajistan/Oman,
# Tayikistan/Haiti, Guayana Francesa/Luxemburgo and Fiyi/Dominica.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-CountryRow($sheet, $row, $country, $stats) {
    # $stats = Casos totales, Nuevos casos, Casos activos, Recuperados,
    #          Casos criticos, Muertes hoy, Muertes  (columns B..H)
    $sheet.Cells.Item($row, 1).Value = $country
    for ($c = 0; $c -lt $stats.Length; $c++) {
        $sheet.Cells.Item($row, 2 + $c).Value = $stats[$c]
    }
}

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 16:52"

# --- Updated statistics (new snapshot), written into their final,
# --- re-sorted row positions -------------------------------------------

# Estados Unidos keeps the #1 spot, just refreshed numbers
Set-CountryRow $ws 4 "Estados Unidos" @(2782539, 2586, 1165206, 1486477, 0, 58, 130856)

# India keeps its position, refreshed numbers
Set-CountryRow $ws 7 "India" @(607344, 2124, 361108, 228363, 0, 25, 17873)

# Portugal's refreshed total (42782) overtakes Kazajistan and Oman,
# so it jumps two places to the top of that trio; Kazajistan and Oman
# each shift down one row keeping their prior (unrefreshed) numbers.
Set-CountryRow $ws 39 "Portugal" @(42782, 328, 28097, 13098, 0, 8, 1587)
Set-CountryRow $ws 40 "Kazajistan" @(42574, 1509, 25533, 16853, 0, 0, 188)
Set-CountryRow $ws 41 "Oman" @(42555, 1361, 25318, 17049, 0, 3, 188)

# Moldavia keeps its position, refreshed numbers
Set-CountryRow $ws 61 "Moldavia" @(16898, 0, 9846, 6501, 0, 2, 551)

# Republica de Macedonia keeps its position, refreshed numbers
Set-CountryRow $ws 81 "Republica de Macedonia" @(6625, 171, 2748, 3556, 0, 15, 321)

# Tayikistan's refreshed total (6058) overtakes Haiti, so it moves one
# row up; Haiti shifts down one row keeping its prior numbers.
Set-CountryRow $ws 84 "Tayikistan" @(6058, 53, 4690, 1316, 0, 0, 52)
Set-CountryRow $ws 85 "Haiti" @(6040, 65, 1032, 4901, 0, 2, 107)

# Guayana Francesa's refreshed total (4444) overtakes Luxemburgo, so it
# moves one row up; Luxemburgo shifts down one row keeping its prior
# numbers.
Set-CountryRow $ws 93 "Guayana Francesa" @(4444, 176, 1680, 2748, 0, 0, 16)
Set-CountryRow $ws 94 "Luxemburgo" @(4345, 0, 4003, 232, 0, 0, 110)

# Cuba keeps its position, refreshed numbers
Set-CountryRow $ws 108 "Cuba" @(2353, 5, 2221, 46, 0, 0, 86)

# Reunion keeps its position, refreshed numbers
Set-CountryRow $ws 155 "Reunion" @(531, 3, 472, 57, 0, 0, 2)

# Mauricio keeps its position, refreshed numbers
Set-CountryRow $ws 159 "Mauricio" @(341, 0, 330, 1, 0, 0, 10)

# Fiyi and Dominica are tied (18/0/18/0/0/0/0 both); the refresh swaps
# their relative order with no numeric change.
Set-CountryRow $ws 205 "Fiyi" @(18, 0, 18, 0, 0, 0, 0)
Set-CountryRow $ws 206 "Dominica" @(18, 0, 18, 0, 0, 0, 0)
